$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain numeric-looking string must keep
# their original "General"/text storage (the source data is scraped text,
# e.g. EU-style "3.320.53" that is not valid as a number at all, alongside
# plain decimals like "253.65" that Excel would otherwise auto-convert to a
# number on assignment). Force text storage via NumberFormat "@" first.

# Row 2
$ws.Range("D2").Value = "97.057.68"
$ws.Range("E2").Value = "  +2.99%  "

# Row 3
$ws.Range("D3").Value = "3.321.73"
$ws.Range("E3").Value = "  +7.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.65"
$ws.Range("E5").Value = "  +8.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.46"
$ws.Range("E6").Value = "  +1.35%  "

# Row 7
$ws.Range("E7").Value = "  +0.79%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.382"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("E9").Value = "  +0.13%  "

# Row 10
$ws.Range("D10").Value = "3.320.53"
$ws.Range("E10").Value = "  +7.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.790"
$ws.Range("E11").Value = "  -4.89%  "

# Row 12
$ws.Range("E12").Value = "  +1.26%  "

# Row 13
$ws.Range("D13").Value = "96.948.51"
$ws.Range("E13").Value = "  +3.13%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.27"
$ws.Range("E14").Value = "  +3.23%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000244"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("E16").Value = "  +6.79%  "

# Row 17
$ws.Range("E17").Value = "  +4.65%  "

# Row 18
$ws.Range("D18").Value = "3.320.97"
$ws.Range("E18").Value = "  +6.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.53"
$ws.Range("E19").Value = "  -3.73%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "480.59"
$ws.Range("E21").Value = "  +7.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000208"
$ws.Range("E22").Value = "  +7.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  -0.50%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  +3.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.59"
$ws.Range("E25").Value = "  +0.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.37"
$ws.Range("E26").Value = "  +2.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.01"
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
$ws.Range("D28").Value = "3.509.84"
$ws.Range("E28").Value = "  +7.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.182"
$ws.Range("E30").Value = "  +1.69%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.237"
$ws.Range("E31").Value = "  -7.11%  "

# Row 32
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$ws.Range("E33").Value = "  -2.58%  "

# Row 34
$ws.Range("E34").Value = "  -1.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.16"
$ws.Range("E35").Value = "  +5.55%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.36"
$ws.Range("E36").Value = "  -5.45%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("E37").Value = "  -4.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "505.84"
$ws.Range("E38").Value = "  +7.88%  "

# Row 39
$ws.Range("E39").Value = "  +2.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.76"
$ws.Range("E40").Value = "  +3.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.447"
$ws.Range("E41").Value = "  -0.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.27"
$ws.Range("E42").Value = "  -0.36%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.28"
$ws.Range("E43").Value = "  +4.32%  "

# Row 44
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.52"
$ws.Range("E44").Value = "  -5.32%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.787"
$ws.Range("E45").Value = "  +15.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.36"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.90"
$ws.Range("E48").Value = "  +1.93%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +5.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.45"
$ws.Range("E50").Value = "  +4.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.48"
$ws.Range("E51").Value = "  +3.33%  "
